{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// 1) Heading: \"Implante de CDI (Evera\u2122 DF4)\" -> \"CDI Dupla C\u00e2mara (Evera DF4)\"\nparagraphs.items[0].insertText(\"CDI Dupla C\u00e2mara (Evera DF4)\", \"Replace\");\n\n// 2) Remove the whole \"Desfibrilador implant\u00e1vel transvenoso DF4.\" paragraph.\nparagraphs.items[1].delete();\n\nawait context.sync();\n\n// Re-load after structural change so indices reflect the new paragraph list.\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// After the deletion, paragraphs are (0) heading, (1) \"Materiais:\",\n// (2) Gerador, (3) Eletrodo de Choque, (4) Eletrodo Atrial, (5) Introdutor.\nparagraphs.items[2].insertText(\"\u2022 Gerador Evera DF4\", \"Replace\");\nparagraphs.items[3].insertText(\"\u2022 Eletrodo 6935M62 DFA\", \"Replace\");\nparagraphs.items[4].insertText(\"\u2022 Eletrodo 5076-52\", \"Replace\");\nparagraphs.items[5].insertText(\"\u2022 Introdutor \u2013 2\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Heading: \"Implante de CDI (Evera\u2122 DF4)\" -> \"CDI Dupla C\u00e2mara (Evera DF4)\"\n$d.Paragraphs(1).Range.Text = \"CDI Dupla C\u00e2mara (Evera DF4)\"\n\n# 2) Remove the whole \"Desfibrilador implant\u00e1vel transvenoso DF4.\" paragraph\n#    (including its paragraph mark) \u2014 it now sits right after the heading.\n$d.Paragraphs(2).Range.Delete()\n\n# After the deletion, paragraphs are:\n#   (1) heading, (2) \"Materiais:\", (3) Gerador, (4) Eletrodo de Choque,\n#   (5) Eletrodo Atrial, (6) Introdutor.\n$d.Paragraphs(3).Range.Text = \"\u2022 Gerador Evera DF4\"\n$d.Paragraphs(4).Range.Text = \"\u2022 Eletrodo 6935M62 DFA\"\n$d.Paragraphs(5).Range.Text = \"\u2022 Eletrodo 5076-52\"\n$d.Paragraphs(6).Range.Text = \"\u2022 Introdutor \u2013 2\"\n"}
